$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds textual price data (already formatted as Text in the source).
# Force Text number format before assigning numeric-looking strings so Excel
# doesn't silently coerce them into numeric cells.
$dCells = @("D2","D4","D5","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D40","D41","D42","D43","D44","D45","D47","D48","D49")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Price updates (column D)
$ws.Range('D2').Value = '248.62'
$ws.Range('D4').Value = '5.628'
$ws.Range('D5').Value = '0.05618'
$ws.Range('D7').Value = '6.464'
$ws.Range('D8').Value = '0.8023'
$ws.Range('D9').Value = '1.064'
$ws.Range('D10').Value = '0.0005740'
$ws.Range('D11').Value = '0.1428'
$ws.Range('D12').Value = '0.07405'
$ws.Range('D13').Value = '0.03199'
$ws.Range('D14').Value = '0.02967'
$ws.Range('D15').Value = '0.09261'
$ws.Range('D16').Value = '0.001669'
$ws.Range('D17').Value = '3.256'
$ws.Range('D18').Value = '0.04693'
$ws.Range('D19').Value = '0.006263'
$ws.Range('D20').Value = '0.001049'
$ws.Range('D21').Value = '0.003812'
$ws.Range('D22').Value = '0.0001500'
$ws.Range('D23').Value = '0.0004600'
$ws.Range('D24').Value = '3.981'
$ws.Range('D25').Value = '2.087'
$ws.Range('D26').Value = '0.3311'
$ws.Range('D27').Value = '0.1277'
$ws.Range('D40').Value = '0.04198'
$ws.Range('D41').Value = '0.006980'
$ws.Range('D42').Value = '0.003500'
$ws.Range('D43').Value = '0.1046'
$ws.Range('D44').Value = '0.009769'
$ws.Range('D45').Value = '0.00005666'
$ws.Range('D47').Value = '0.6801'
$ws.Range('D48').Value = '0.02916'
$ws.Range('D49').Value = '0.00002100'

# Text field updates (Coin name / Link / Volume label columns)
$ws.Range('B10').Value = 'One'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('E10').Value = '9OneONE'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('E11').Value = '10WazirXWRX'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'
$ws.Range('B13').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C13').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('E13').Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B19').Value = 'TigerCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('E19').Value = '18TigerCashTCH'
$ws.Range('B20').Value = 'BitKan'
$ws.Range('C20').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('E20').Value = '19BitKanKAN'
$ws.Range('B21').Value = 'HotbitToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('E21').Value = '20HotbitTokenHTB'
$ws.Range('B22').Value = 'NitroEx'
$ws.Range('C22').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('E22').Value = '21NitroExNTX'
$ws.Range('B23').Value = 'UpBots'
$ws.Range('C23').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('E23').Value = '22UpBotsUBXT'
$ws.Range('B24').Value = 'LEO'
$ws.Range('C24').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('E24').Value = '23LEOLEO'
$ws.Range('B25').Value = 'BTSEToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('E25').Value = '24BTSETokenBTSE'
$ws.Range('B26').Value = 'BitpandaEcosystemToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('E26').Value = '25BitpandaEcosystemTokenBEST'
$ws.Range('B27').Value = 'ProBitToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('E27').Value = '26ProBitTokenPROB'
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('E42').Value = '41CEJICEJI'
$ws.Range('B43').Value = 'BKEXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('E43').Value = '42BKEXTokenBKK'
$ws.Range('E48').Value = '47BOLOBOLOWorstin24h'
